# Fixed update to excel issue
# - Rename "Requested quantity" headers on the existing sheets
# - Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

# --- Update header labels on existing sheets -------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet ----------------------------------------
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$data = @(
    @(45039.99999999999, 38, 4.724196028557685, 69.48386525379559),
    @(45046.99999999999, 38, 6.036403752373369, 72.15069043208655),
    @(45053.99999999999, 38, 4.344411185869462, 71.66746680044032),
    @(45060.99999999999, 37, 4.468457462162412, 72.69763229161141),
    @(45074.99999999999, 37, 5.985447843216842, 67.91282199265348),
    @(45081.99999999999, 37, 3.26497498227738, 68.19584582704151),
    @(45088.99999999999, 36, 1.238636608320768, 69.52021560521634),
    @(45095.99999999999, 36, 2.464733055337547, 69.61010310956354),
    @(45102.99999999999, 36, 2.990201137511742, 71.23684861118853),
    @(45116.99999999999, 35, 1.707703859705789, 66.91591645686914),
    @(45123.99999999999, 35, 1.621680180779342, 67.87417722941368),
    @(45144.99999999999, 34, -1.209347626187893, 68.28362214255951),
    @(45151.99999999999, 33, 0.6204628885127841, 68.8996786922621),
    @(45158.99999999999, 33, -2.942264566244702, 64.85507197258316),
    @(45165.99999999999, 33, -2.26290456718343, 67.81507885668363),
    @(45186.99999999999, 32, -1.213508478815205, 65.95849725667001),
    @(45193.99999999999, 32, 0.04549005304805519, 66.89543271648124),
    @(45200.99999999999, 31, -0.510133002687566, 65.61902817448012),
    @(45207.99999999999, 31, -2.862120874982139, 64.84032189019071),
    @(45214.99999999999, 31, -2.026367455784517, 65.0435080612598),
    @(45221.99999999999, 30, -5.454180066051711, 64.33345426202355),
    @(45228.99999999999, 30, -2.905340922178822, 66.60254390638971),
    @(45235.99999999999, 30, -3.304976169951232, 63.50043699026571),
    @(45242.99999999999, 29, -5.146597635960095, 62.4651033347796),
    @(45249.99999999999, 29, -4.925004572500977, 63.66344284708509),
    @(45256.99999999999, 29, -6.14108205761193, 61.08516251968904)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsForecast.Range("A2:A27").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Move the new sheet to the end so the tab order matches: Weekly Quantity,
# Monthly Trend, PO Forecast
$wsForecast.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
